# "LAST TIME for now"
# Adds a new person record (Shah Rukh / 20230319SH) as row 7 on Sheet1,
# re-zooms/re-selects the sheet view, and touches the page setup
# (orientation) so the stale printer-settings relationship is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Append the new data row ------------------------------------------------
$ws.Range("A7").Value = "20230319SH"
$ws.Range("B7").Value = "Shah Rukh"
$ws.Range("C7").Value = "jenishkanani93@gmail.com"
$ws.Range("D7").Value = 8866432894
$ws.Range("E7").Value = "abc, st bus stand"
$ws.Range("F7").Value = 1234
$ws.Range("G7").Value = 12414
$ws.Range("H7").Value = "GSRTC"
$ws.Range("I7").Value = 123331

# --- View changes: zoom to 140% and move the selection to A5 ---------------
$excel.ActiveWindow.Zoom = 140
$ws.Range("A5").Select()

# --- Page setup: re-apply orientation (drops the old printer-settings ref) -
$ws.PageSetup.Orientation = 1
